$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking "cryptos" snapshot refresh (GitHub Actions bot).
# Price (column D) cells are stored as literal text (e.g. "1.630.19" is
# thousand-dot-separated, not a real number), so force a Text number
# format before assigning to stop Excel from re-interpreting strings
# such as "240.80" or "1.00" as numeric values and trimming digits.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.904.60"
$ws.Range("E2").Value = "  +0.68%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.631.37"
$ws.Range("E3").Value = "  +1.05%  "

$ws.Range("E4").Value = "  +0.62%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.66"
$ws.Range("E5").Value = "  +0.86%  "

$ws.Range("E6").Value = "  -0.72%  "

$ws.Range("E7").Value = "  +0.64%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "28.61"
$ws.Range("E8").Value = "  -1.10%  "

$ws.Range("E9").Value = "  -0.26%  "

$ws.Range("E10").Value = "  -0.16%  "

$ws.Range("E11").Value = "  -0.35%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.864.75"
$ws.Range("E12").Value = "  +1.19%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.629.09"
$ws.Range("E13").Value = "  +0.86%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.562"
$ws.Range("E14").Value = "  -0.69%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "9.26"
$ws.Range("E15").Value = "  +7.03%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.906.15"
$ws.Range("E16").Value = "  +0.61%  "

$ws.Range("E17").Value = "  -0.33%  "

$ws.Range("E18").Value = "  -1.23%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "240.80"
$ws.Range("E19").Value = "  -0.44%  "

$ws.Range("E20").Value = "  -0.85%  "

$ws.Range("E21").Value = "  +0.51%  "

$ws.Range("E22").Value = "  +0.83%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.79"
$ws.Range("E23").Value = "  +1.54%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.15"
$ws.Range("E24").Value = "  +1.34%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.83"
$ws.Range("E25").Value = "  +0.81%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.44"
$ws.Range("E26").Value = "  -1.23%  "

$ws.Range("E27").Value = "  -0.62%  "

$ws.Range("E28").Value = "  +0.00%  "

$ws.Range("E29").Value = "  +0.53%  "

$ws.Range("E30").Value = "  +1.24%  "

$ws.Range("E31").Value = "  +3.09%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.37"
$ws.Range("E32").Value = "  +2.09%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.17"
$ws.Range("E33").Value = "  -0.31%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.424.07"
$ws.Range("E34").Value = "  -0.43%  "

$ws.Range("E35").Value = "  +2.98%  "

$ws.Range("E36").Value = "  -1.55%  "

$ws.Range("E37").Value = "  -4.81%  "

$ws.Range("E38").Value = "  -0.31%  "

$ws.Range("E39").Value = "  -0.19%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "75.16"
$ws.Range("E40").Value = "  +8.05%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.554"
$ws.Range("E41").Value = "  -0.54%  "

$ws.Range("E42").Value = "  +0.05%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0499"
$ws.Range("E43").Value = "  -0.99%  "

$ws.Range("E44").Value = "  -0.18%  "

$ws.Range("E45").Value = "  +0.64%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("E46").Value = "  +0.01%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "49.99"
$ws.Range("E47").Value = "  -7.41%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.772.68"
$ws.Range("E48").Value = "  +1.24%  "

$ws.Range("E49").Value = "  -2.38%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "90.46"
$ws.Range("E50").Value = "  +3.62%  "

$ws.Range("E51").Value = "  +7.09%  "
